$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.02328103994796504
$ws.Range("B2").Value = -0.06201132485515971

$ws.Range("A3").Value = 0.05483536340375555
$ws.Range("B3").Value = 0.07646821791777769

$ws.Range("A4").Value = 0.01644714489527688
$ws.Range("B4").Value = -0.04297431422540074

$ws.Range("A5").Value = 0.1601979657762684
$ws.Range("B5").Value = 0.1753738058469409

$ws.Range("A6").Value = -0.04883717205316224
$ws.Range("B6").Value = -0.04483938565520964

$ws.Range("A7").Value = 0.5011471247107649
$ws.Range("B7").Value = 0.492171605091032

$ws.Range("A8").Value = 0.2440913612659691
$ws.Range("B8").Value = 0.2559863711442495
